# "New building shells in, focus on Implementation next"
#
# Rotates the Void-building columns (E: Faction Embassy, F: Abyssal
# Pathfinder, G: Void Rudder) on the "Specialty Buildings" sheet one step
# to the left (E<-F, F<-G, G<-E) for every populated row, then leaves the
# "Specialty Buildings" sheet as the active tab/selection, matching where
# the author's focus moved to next.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specialty Buildings")

$rows = @(1, 6, 7, 8, 9, 10, 11, 12, 13)

foreach ($r in $rows) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2

    $eCell.Value = $fVal
    $fCell.Value = $gVal
    $gCell.Value = $eVal
}

# Move focus to the "Specialty Buildings" sheet, landing on column H
# (Weather Manipulator), which is where the next round of work continues.
$ws.Activate()
$ws.Range("H1:H1048576").Select()
